$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5499.7
$ws.Range("I86").Value = 4999.3335
$ws.Range("J86").Value = 5714.143
$ws.Range("K86").Value = 4999.3335
$ws.Range("L86").Value = 5714.143
$ws.Range("M86").Value = -3876.3335
$ws.Range("N86").Value = -7960.143

$ws.Range("H89").Value = 5499.7
$ws.Range("I89").Value = 4999.3335
$ws.Range("J89").Value = 5714.143
$ws.Range("K89").Value = 24996.6675
$ws.Range("L89").Value = 28570.715
$ws.Range("M89").Value = -19380.6675
$ws.Range("N89").Value = -39802.715

$ws.Range("H125").Value = 3233
$ws.Range("I125").Value = 2966
$ws.Range("K125").Value = 26694
$ws.Range("M125").Value = -24234

$ws.Range("H132").Value = 3028.2258
$ws.Range("I132").Value = 1874.7931
$ws.Range("K132").Value = 5624.379300000001
$ws.Range("M132").Value = -3094.379300000001

$ws.Range("H137").Value = 3417.6
$ws.Range("I137").Value = 2697.5334
$ws.Range("J137").Value = 4857.7334
$ws.Range("K137").Value = 8092.600199999999
$ws.Range("L137").Value = 14573.2002
$ws.Range("M137").Value = -5542.600199999999
$ws.Range("N137").Value = -19673.2002

$ws.Range("H138").Value = 3559.963
$ws.Range("J138").Value = 3878.6667
$ws.Range("L138").Value = 11636.0001
$ws.Range("N138").Value = -21916.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2784.2856
$ws.Range("I26").Value = 294.8
$ws.Range("K26").Value = 294.8
$ws.Range("M26").Value = 35.19999999999999

$ws.Range("H43").Value = 29067.25
$ws.Range("I43").Value = 27999
$ws.Range("J43").Value = 29423.334
$ws.Range("K43").Value = 27999
$ws.Range("L43").Value = 29423.334
$ws.Range("M43").Value = -27686
$ws.Range("N43").Value = -30049.334

$ws.Range("H45").Value = 52634940
$ws.Range("I45").Value = 66668816
$ws.Range("J45").Value = 7903.5
$ws.Range("K45").Value = 66668816
$ws.Range("L45").Value = 7903.5
$ws.Range("M45").Value = -66668439
$ws.Range("N45").Value = -8657.5

$ws.Range("H50").Value = 4577.5
$ws.Range("J50").Value = 4577.5
$ws.Range("L50").Value = 4577.5
$ws.Range("N50").Value = -6005.5

$ws.Range("H102").Value = 1858.0555
$ws.Range("I102").Value = 1858.0555
$ws.Range("K102").Value = 1858.0555
$ws.Range("M102").Value = -236.0554999999999

$ws.Range("H132").Value = 6230.457
$ws.Range("I132").Value = 4547.3667
$ws.Range("K132").Value = 13642.1001
$ws.Range("M132").Value = -11112.1001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 55000
$ws.Range("J13").Value = 55000
$ws.Range("L13").Value = 55000
$ws.Range("N13").Value = -55336

$ws.Range("H86").Value = 2090.0557
$ws.Range("I86").Value = 965.2857
$ws.Range("J86").Value = 6026.75
$ws.Range("K86").Value = 965.2857
$ws.Range("L86").Value = 6026.75
$ws.Range("M86").Value = 157.7143
$ws.Range("N86").Value = -8272.75

$ws.Range("H89").Value = 2090.0557
$ws.Range("I89").Value = 965.2857
$ws.Range("J89").Value = 6026.75
$ws.Range("K89").Value = 4826.4285
$ws.Range("L89").Value = 30133.75
$ws.Range("M89").Value = 789.5715
$ws.Range("N89").Value = -41365.75

$ws.Range("H94").Value = 899.2
$ws.Range("I94").Value = 919.3
$ws.Range("J94").Value = 859
$ws.Range("K94").Value = 919.3
$ws.Range("L94").Value = 859
$ws.Range("M94").Value = -468.3
$ws.Range("N94").Value = -1761

$ws.Range("H99").Value = 5524.25
$ws.Range("I99").Value = 6479
$ws.Range("K99").Value = 6479
$ws.Range("M99").Value = -4981

$ws.Range("H105").Value = 23136.25
$ws.Range("I105").Value = 27899.75
$ws.Range("K105").Value = 27899.75
$ws.Range("M105").Value = -26152.75

$ws.Range("H134").Value = 2467.625
$ws.Range("I134").Value = 1498.2142
$ws.Range("K134").Value = 4494.642599999999
$ws.Range("M134").Value = -1959.642599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26513.375
$ws.Range("I31").Value = 3779.818
$ws.Range("J31").Value = 33272
$ws.Range("K31").Value = 3779.818
$ws.Range("L31").Value = 33272
$ws.Range("M31").Value = -3484.818
$ws.Range("N31").Value = -33862

$ws.Range("H34").Value = 26513.375
$ws.Range("I34").Value = 3779.818
$ws.Range("J34").Value = 33272
$ws.Range("K34").Value = 3779.818
$ws.Range("L34").Value = 33272
$ws.Range("M34").Value = -3577.818
$ws.Range("N34").Value = -33676

$ws.Range("H86").Value = 10630.4
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 12913
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 12913
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -15159

$ws.Range("H89").Value = 10630.4
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 12913
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 64565
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -75797

$ws.Range("H132").Value = 3829.111
$ws.Range("I132").Value = 3105.1667
$ws.Range("J132").Value = 5277
$ws.Range("K132").Value = 9315.500100000001
$ws.Range("L132").Value = 15831
$ws.Range("M132").Value = -6785.500100000001
$ws.Range("N132").Value = -20891

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8329.166999999999
$ws.Range("I70").Value = 7499
$ws.Range("J70").Value = 9159.333000000001
$ws.Range("K70").Value = 7499
$ws.Range("L70").Value = 9159.333000000001
$ws.Range("M70").Value = -7229
$ws.Range("N70").Value = -9699.333000000001

$ws.Range("H73").Value = 8329.166999999999
$ws.Range("I73").Value = 7499
$ws.Range("J73").Value = 9159.333000000001
$ws.Range("K73").Value = 7499
$ws.Range("L73").Value = 9159.333000000001
$ws.Range("M73").Value = -6563
$ws.Range("N73").Value = -11031.333

$ws.Range("H80").Value = 5868.357
$ws.Range("I80").Value = 4069.625
$ws.Range("J80").Value = 8266.666999999999
$ws.Range("K80").Value = 4069.625
$ws.Range("L80").Value = 8266.666999999999
$ws.Range("M80").Value = -3071.625
$ws.Range("N80").Value = -10262.667

$ws.Range("H83").Value = 5868.357
$ws.Range("I83").Value = 4069.625
$ws.Range("J83").Value = 8266.666999999999
$ws.Range("K83").Value = 20348.125
$ws.Range("L83").Value = 41333.335
$ws.Range("M83").Value = -15356.125
$ws.Range("N83").Value = -51317.335

$ws.Range("H113").Value = 6375.6665
$ws.Range("I113").Value = 1918
$ws.Range("K113").Value = 1918
$ws.Range("M113").Value = 252

$ws.Range("H122").Value = 12495.947
$ws.Range("I122").Value = 16378.777
$ws.Range("J122").Value = 9001.4
$ws.Range("K122").Value = 49136.331
$ws.Range("L122").Value = 27004.2
$ws.Range("M122").Value = -46686.331
$ws.Range("N122").Value = -31904.2

$ws.Range("H132").Value = 5103.9375
$ws.Range("I132").Value = 4263.8
$ws.Range("J132").Value = 6504.1665
$ws.Range("K132").Value = 12791.4
$ws.Range("L132").Value = 19512.4995
$ws.Range("M132").Value = -10261.4
$ws.Range("N132").Value = -24572.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10534.5
$ws.Range("I7").Value = 6569.7334
$ws.Range("K7").Value = 6569.7334
$ws.Range("M7").Value = -6457.7334

$ws.Range("H46").Value = 3855.8572
$ws.Range("I46").Value = 1640
$ws.Range("J46").Value = 6071.7144
$ws.Range("K46").Value = 1640
$ws.Range("L46").Value = 6071.7144
$ws.Range("M46").Value = -1452
$ws.Range("N46").Value = -6447.7144

$ws.Range("H82").Value = 5676.1113
$ws.Range("I82").Value = 5146.375
$ws.Range("J82").Value = 6099.9
$ws.Range("K82").Value = 5146.375
$ws.Range("L82").Value = 6099.9
$ws.Range("M82").Value = -4785.375
$ws.Range("N82").Value = -6821.9

$ws.Range("H85").Value = 5676.1113
$ws.Range("I85").Value = 5146.375
$ws.Range("J85").Value = 6099.9
$ws.Range("K85").Value = 5146.375
$ws.Range("L85").Value = 6099.9
$ws.Range("M85").Value = -3898.375
$ws.Range("N85").Value = -8595.9

$ws.Range("H93").Value = 13949.861
$ws.Range("I93").Value = 11357.091
$ws.Range("J93").Value = 18024.215
$ws.Range("K93").Value = 11357.091
$ws.Range("L93").Value = 18024.215
$ws.Range("M93").Value = -10109.091
$ws.Range("N93").Value = -20520.215

$ws.Range("H100").Value = 5132.6
$ws.Range("I100").Value = 895
$ws.Range("J100").Value = 7957.6665
$ws.Range("K100").Value = 895
$ws.Range("L100").Value = 7957.6665
$ws.Range("M100").Value = -354
$ws.Range("N100").Value = -9039.666499999999

$ws.Range("H126").Value = 10534.5
$ws.Range("I126").Value = 6569.7334
$ws.Range("K126").Value = 19709.2002
$ws.Range("M126").Value = -17239.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1288.4
$ws.Range("J100").Value = 1801.5
$ws.Range("L100").Value = 3603
$ws.Range("N100").Value = -4685

$ws.Range("H122").Value = 2796.0688
$ws.Range("I122").Value = 1528.625
$ws.Range("J122").Value = 8879.799999999999
$ws.Range("K122").Value = 4585.875
$ws.Range("L122").Value = 26639.4
$ws.Range("M122").Value = -2135.875
$ws.Range("N122").Value = -31539.4

$ws.Range("H126").Value = 4206.8696
$ws.Range("I126").Value = 3610
$ws.Range("J126").Value = 4525.2
$ws.Range("K126").Value = 10830
$ws.Range("L126").Value = 13575.6
$ws.Range("M126").Value = -8360
$ws.Range("N126").Value = -18515.6

$ws.Range("H136").Value = 4202.3477
$ws.Range("I136").Value = 2127.1538
$ws.Range("J136").Value = 6900.1
$ws.Range("K136").Value = 6381.4614
$ws.Range("L136").Value = 20700.3
$ws.Range("M136").Value = -3831.4614
$ws.Range("N136").Value = -25800.3
